# Phân công lại cho mọi người nhé
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time slots for rows 7-14 to "7h-8h"
$ws.Range("H7").Value = "7h-8h"
$ws.Range("H8").Value = "7h-8h"
$ws.Range("H9").Value = "7h-8h"
$ws.Range("H10").Value = "7h-8h"
$ws.Range("H11").Value = "7h-8h"
$ws.Range("H12").Value = "7h-8h"
$ws.Range("H13").Value = "7h-8h"
$ws.Range("H14").Value = "7h-8h"

# Row 15: reassign person from "Ngọc" to "Lan" and time from "10h-11h" to "8h-9h"
$ws.Range("G15").Value = "Lan"
$ws.Range("H15").Value = "8h-9h"

# Update the current selection to reflect where the edit ended (F15)
$ws.Range("F15").Select()
